# Working-hours log: "fixed bug in detection of convergence"
# The previous cut-off (row 144) was missing its clock-out time, and the
# following real shift (now row 145) had been left out of the sheet
# entirely, sitting in what used to be a blank spacer row. Filling these
# in pushes the blank spacer row and the three summary rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the missing entry: insert a new row at 145, shifting the
# blank spacer row and the summary rows (old 145-148) down to 146-149.
$null = $ws.Rows("145:145").Insert()

# Complete row 144 (clock-out time) and add row 145 (the shift that had
# been missing from the log entirely).
$ws.Range("E144").Value = 0.79166666666666663
$ws.Range("A145").Value = 2014
$ws.Range("B145").Value = 7
$ws.Range("C145").Value = 21
$ws.Range("D145").Value = 0.83333333333333337
$ws.Range("E145").Value = 0.875

# Extend the "time spent" / "sum [h]" formulas down through the two rows
# (mirrors dragging the fill handle from F143:G143 down to F145:G145).
$ws.Range("F144:F145").Formula = "=(E144-D144)*24*60"
$ws.Range("G144:G145").Formula = "=F144/60"

# Summary block (now rows 147-149) needs its SUM range extended to
# include the newly-added row.
$ws.Range("F147").Formula = "=SUM(F2:F145)"

# Restore the view to roughly where it was (scrolled near the bottom,
# with the previously-edited cell selected).
$excel.ActiveWindow.ScrollRow = 127
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("F145").Select()
